# Fixed gn and problem for others
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: Length of the solution, Length of the search path, Execution time (seconds)
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 152.9034914

# Remove rows 3 and 4 (the other puzzles), shrinking the used range to A1:F2
$ws.Range("A3:F4").Delete()
